# Edit: B1--and-B2-PowerPoint.pptx
#
# 1) Slide 5's table switches to a different table style (tableStyleId).
# 2) The deck's theme (color scheme) used by the slide master changes from
#    the "Integral" / "Red Violet" palette to the standard "Office" palette.
#
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 -------------------------------------------
# slide5.xml is the 5th slide in the deck; its table (graphicFrame) is the
# 2nd shape on that slide (shape 1 is the title textbox, shape 3 is another
# textbox).
$s = $p.Slides.Item(5)
$tblShape = $s.Shapes.Item(2)
if ($tblShape.HasTable) {
    $tblShape.Table.ApplyStyle("{D0C113BA-3858-4C08-9027-53834E927499}")
}

# --- 2) Swap the theme colour palette applied to the deck -----------------
# Replace the current "Red Violet" theme colours (used by the slide master /
# theme1.xml) with the standard Office theme colours. Colors(1..12) map to
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink in that order.
function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

Set-ThemeColor $colorScheme 1  0   0   0     # dk1      000000
Set-ThemeColor $colorScheme 2  255 255 255   # lt1      FFFFFF
Set-ThemeColor $colorScheme 3  68  84  106   # dk2      44546A
Set-ThemeColor $colorScheme 4  231 230 230   # lt2      E7E6E6
Set-ThemeColor $colorScheme 5  91  155 213   # accent1  5B9BD5
Set-ThemeColor $colorScheme 6  237 125 49    # accent2  ED7D31
Set-ThemeColor $colorScheme 7  165 165 165   # accent3  A5A5A5
Set-ThemeColor $colorScheme 8  255 192 0     # accent4  FFC000
Set-ThemeColor $colorScheme 9  68  114 196   # accent5  4472C4
Set-ThemeColor $colorScheme 10 112 173 71    # accent6  70AD47
Set-ThemeColor $colorScheme 11 5   99  193   # hlink    0563C1
Set-ThemeColor $colorScheme 12 149 79  114   # folHlink 954F72
